$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# Merge the two runs "...more p" / "rocessors are added. " back into a
# single run and drop the _GoBack bookmark that used to sit between them.
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("more processors are added", $true, $false, $false, $false, $false, $true, 1, $false, "more processors are added", 2)

# --- Change 2 ---------------------------------------------------------
# "Browse and search the full collection at" ->
#   "Browse and search the full " + (italic run) "curriculum" + (italic run) " at"
$outer2 = $d.Content
$found2 = $outer2.Find.Execute("full collection at")
$word2 = $d.Range($outer2.Start, $outer2.End)
$foundWord2 = $word2.Find.Execute("collection")
$word2.Text = "curriculum"

# Force the "curriculum" word onto its own run by toggling an explicit
# character formatting change (off then back on) - this prevents the host
# from silently re-coalescing it with the identically-formatted text
# around it.
$outer2b = $d.Content
$found2b = $outer2b.Find.Execute("curriculum at")
$cur2 = $d.Range($outer2b.Start, $outer2b.Start + 10)
$cur2.Italic = 0
$cur2.Italic = 1

# Likewise keep the following <w:br/> as its own run (it sits right after
# " at" in the original markup and must stay untouched/separate).
$outer2c = $d.Content
$found2c = $outer2c.Find.Execute("curriculum at")
$br2 = $d.Range($outer2c.End, $outer2c.End + 1)
$br2.Italic = 0
$br2.Italic = 1

# --- Change 3 ---------------------------------------------------------
# "material and the rest of the collection in our GitHub repository at" ->
#   "material and the rest of the " + (italic run) "curriculum" +
#   _GoBack bookmark + (italic run) " in our GitHub repository at"
$outer3 = $d.Content
$found3 = $outer3.Find.Execute("rest of the collection in our GitHub")
$word3 = $d.Range($outer3.Start, $outer3.End)
$foundWord3 = $word3.Find.Execute("collection")
$word3.Text = "curriculum"

$outer3b = $d.Content
$found3b = $outer3b.Find.Execute("curriculum in our GitHub repository at")
$cur3 = $d.Range($outer3b.Start, $outer3b.Start + 10)
$cur3.Italic = 0
$cur3.Italic = 1

$outer3c = $d.Content
$found3c = $outer3c.Find.Execute("curriculum in our GitHub repository at")
$tail3 = $d.Range($outer3c.Start + 10, $outer3c.End)
$tail3.Italic = 0
$tail3.Italic = 1

# Re-insert the _GoBack bookmark right after the newly split "curriculum"
# run (this is where it lives in the target document).
$outer3d = $d.Content
$found3d = $outer3d.Find.Execute("curriculum in our GitHub repository at")
$bmPoint = $d.Range($outer3d.Start + 10, $outer3d.Start + 10)
$d.Bookmarks.Add("_GoBack", $bmPoint)
